$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text change (Q1)
$ws.Range("Q1").Value = "wtkappa.scale_trim"

# Row 2 (All data)
$ws.Range("B2").Value = [double]"-8.43769498715119e-17"
$ws.Range("C2").Value = [double]"-1.576516694967722e-16"
$ws.Range("Q2").Value = 0.7808705382933534

# Row 3 (QUESTION_1)
$ws.Range("B3").Value = -0.01149607432881209
$ws.Range("C3").Value = -0.01543341287607883
$ws.Range("Q3").Value = 0.7896756462802278

# Row 4 (QUESTION_2)
$ws.Range("B4").Value = 0.006754119518834534
$ws.Range("C4").Value = 0.01028894191738561
$ws.Range("Q4").Value = 0.7665395469417655

# Row 5 (QUESTION_3)
$ws.Range("B5").Value = 0.005058824988670463
$ws.Range("C5").Value = -0.04115576766954331
$ws.Range("Q5").Value = 0.7811601973293526

# Row 6 (QUESTION_4)
$ws.Range("B6").Value = -0.001108556945570738
$ws.Range("C6").Value = -0.04115576766954331
$ws.Range("Q6").Value = 0.8093308458669355

# Row 7 (QUESTION_5)
$ws.Range("B7").Value = 0.0007916867668774308
$ws.Range("C7").Value = 0.08745600629777898
$ws.Range("Q7").Value = 0.7793853523195861
